$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 512.25
$ws.Range("I58").Value = 512.25
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1536.75
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1386.75
$ws.Range("H103").Value = 2164.889
$ws.Range("I103").Value = 966.6667
$ws.Range("J103").Value = 2764
$ws.Range("K103").Value = 2900.0001
$ws.Range("L103").Value = 8292
$ws.Range("M103").Value = -2314.0001
$ws.Range("H113").Value = 2555.4443
$ws.Range("I113").Value = 3850
$ws.Range("J113").Value = 2185.5715
$ws.Range("K113").Value = 3850
$ws.Range("L113").Value = 2185.5715
$ws.Range("M113").Value = -596
$ws.Range("N113").Value = -8693.5715
$ws.Range("H127").Value = 6115
$ws.Range("I127").Value = 6115
$ws.Range("J127").Value = 0
$ws.Range("K127").Value = 18345
$ws.Range("L127").Value = 0
$ws.Range("M127").Value = -13385
$ws.Range("H137").Value = 1958.963
$ws.Range("I137").Value = 1463.2727
$ws.Range("J137").Value = 4140
$ws.Range("K137").Value = 4389.8181
$ws.Range("L137").Value = 12420
$ws.Range("M137").Value = -1839.8181
$ws.Range("N137").Value = -17520
$ws.Range("H141").Value = 2258.5
$ws.Range("I141").Value = 2258.5
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 6775.5
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -1595.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1901.375
$ws.Range("I102").Value = 1173
$ws.Range("J102").Value = 7000
$ws.Range("K102").Value = 1173
$ws.Range("L102").Value = 7000
$ws.Range("M102").Value = 449
$ws.Range("H140").Value = 56666.332
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 56666.332
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 56666.332
$ws.Range("N140").Value = -67026.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 4253.643
$ws.Range("I94").Value = 3955.1
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 3955.1
$ws.Range("L94").Value = 5000
$ws.Range("M94").Value = -3504.1
$ws.Range("H105").Value = 2221.1333
$ws.Range("I105").Value = 1678.2307
$ws.Range("J105").Value = 5750
$ws.Range("K105").Value = 1678.2307
$ws.Range("L105").Value = 5750
$ws.Range("M105").Value = 68.76929999999993

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 10750
$ws.Range("I29").Value = 6500
$ws.Range("J29").Value = 15000
$ws.Range("K29").Value = 6500
$ws.Range("L29").Value = 15000
$ws.Range("M29").Value = -6207
$ws.Range("N29").Value = -15586
$ws.Range("H41").Value = 25333.25
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 25333.25
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 25333.25
$ws.Range("N41").Value = -26189.25
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H94").Value = 2559
$ws.Range("I94").Value = 2183.6667
$ws.Range("J94").Value = 3403.5
$ws.Range("K94").Value = 2183.6667
$ws.Range("L94").Value = 3403.5
$ws.Range("M94").Value = -1732.6667
$ws.Range("N94").Value = -4305.5
$ws.Range("H105").Value = 895.6
$ws.Range("I105").Value = 870.25
$ws.Range("J105").Value = 997
$ws.Range("K105").Value = 870.25
$ws.Range("L105").Value = 997
$ws.Range("M105").Value = 876.75
$ws.Range("N105").Value = -4491
$ws.Range("H132").Value = 1572.4
$ws.Range("I132").Value = 1758.2413
$ws.Range("J132").Value = 674.1667
$ws.Range("K132").Value = 5274.7239
$ws.Range("L132").Value = 2022.5001
$ws.Range("M132").Value = -2744.7239
$ws.Range("N132").Value = -7082.5001
$ws.Range("H134").Value = 1269.5714
$ws.Range("I134").Value = 1269.5714
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3808.7142
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1273.7142

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 375000500
$ws.Range("I4").Value = 333333860
$ws.Range("J4").Value = 500000400
$ws.Range("K4").Value = 1000001580
$ws.Range("L4").Value = 1500001200
$ws.Range("M4").Value = -1000001468
$ws.Range("H128").Value = 288933
$ws.Range("I128").Value = 288933
$ws.Range("J128").Value = 0
$ws.Range("K128").Value = 866799
$ws.Range("L128").Value = 0
$ws.Range("M128").Value = -861819

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 713.2857
$ws.Range("I2").Value = 1626.4445
$ws.Range("J2").Value = 280.73685
$ws.Range("K2").Value = 1626.4445
$ws.Range("L2").Value = 280.73685
$ws.Range("M2").Value = -1513.4445
$ws.Range("N2").Value = -506.73685
$ws.Range("H70").Value = 10002440
$ws.Range("I70").Value = 12502387
$ws.Range("J70").Value = 2650
$ws.Range("K70").Value = 12502387
$ws.Range("L70").Value = 2650
$ws.Range("M70").Value = -12502117
$ws.Range("N70").Value = -3190
$ws.Range("H73").Value = 10002440
$ws.Range("I73").Value = 12502387
$ws.Range("J73").Value = 2650
$ws.Range("K73").Value = 12502387
$ws.Range("L73").Value = 2650
$ws.Range("M73").Value = -12501451
$ws.Range("N73").Value = -4522

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1689.1111
$ws.Range("I46").Value = 1364.7059
$ws.Range("J46").Value = 2240.6
$ws.Range("K46").Value = 1364.7059
$ws.Range("L46").Value = 2240.6
$ws.Range("M46").Value = -1176.7059
$ws.Range("H62").Value = 23750
$ws.Range("I62").Value = 10000
$ws.Range("J62").Value = 37500
$ws.Range("K62").Value = 10000
$ws.Range("L62").Value = 37500
$ws.Range("M62").Value = -9376
$ws.Range("H65").Value = 23750
$ws.Range("I65").Value = 10000
$ws.Range("J65").Value = 37500
$ws.Range("K65").Value = 30000
$ws.Range("L65").Value = 112500
$ws.Range("M65").Value = -26880
$ws.Range("H82").Value = 1404
$ws.Range("H85").Value = 1404
$ws.Range("H122").Value = 5541.615
$ws.Range("I122").Value = 5397.5
$ws.Range("J122").Value = 6022
$ws.Range("K122").Value = 16192.5
$ws.Range("L122").Value = 18066
$ws.Range("M122").Value = -13742.5
$ws.Range("H132").Value = 4491
$ws.Range("I132").Value = 4622.6
$ws.Range("J132").Value = 3833
$ws.Range("K132").Value = 13867.8
$ws.Range("L132").Value = 11499
$ws.Range("M132").Value = -11337.8
$ws.Range("N132").Value = -16559
$ws.Range("H134").Value = 82343.2
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 82343.2
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 82343.2
$ws.Range("N134").Value = -92483.2
$ws.Range("H138").Value = 79801.836
$ws.Range("I138").Value = 80000
$ws.Range("J138").Value = 78811
$ws.Range("K138").Value = 80000
$ws.Range("L138").Value = 78811
$ws.Range("M138").Value = -74860
$ws.Range("N138").Value = -89091

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 299.63635
$ws.Range("I107").Value = 329
$ws.Range("J107").Value = 248.25
$ws.Range("K107").Value = 987
$ws.Range("L107").Value = 744.75
$ws.Range("M107").Value = 933
$ws.Range("N107").Value = -4584.75
$ws.Range("H132").Value = 1354.7273
$ws.Range("I132").Value = 990.2
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 2970.6
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -440.6000000000004
